## Timesheet update: add a missed time entry on "Week 6" (row 6) and let the
## existing running-total formulas (E20/E21 on every week sheet, cascading
## week to week) recompute automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Week 6")
$ws.Activate()

# Row 4 already carries the formatting the new row needs (date format on A,
# time format on B/C, wrap-text description column on D, ht=26 for the
# two-line description) - clone it onto row 6 before filling in the values.
$ws.Range("A4:E4").Copy()
$ws.Range("A6:E6").PasteSpecial(-4122)  # xlPasteFormats
$ws.Rows.Item(6).RowHeight = 26

$ws.Range("A6").Value = 41684
$ws.Range("B6").Value = 0.41666666666666669
$ws.Range("C6").Value = 0.5
$ws.Range("D6").Value = "Uploaded product images, created states lookup table, modified table relationships, resolved terminal issue"
$ws.Range("E6").Value = 2

# Matches the author's final selection in the saved file.
$ws.Range("D6").Select()
